$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = "[53.6362303017758, 71.85277155970623]"
$ws.Range("U2").Value = "[42.45159382969866, 54.86822623054531]"

# Row 3
$ws.Range("M3").Value = "[50.872209087677334, 74.23044147044797]"
$ws.Range("N3").Value = [double]"4.596323321948148e-14"
$ws.Range("O3").Value = [double]"4.596323321948148e-14"
$ws.Range("Q3").Value = "[1.0629212381515023, 1.4402897250691948]"
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0
$ws.Range("U3").Value = "[41.12999394560171, 54.31935586222566]"
$ws.Range("Y3").Value = 19.56986986987028
$ws.Range("Z3").Value = 21.09479479479523

# Row 4
$ws.Range("M4").Value = "[51.10217547189261, 76.32384407869816]"
$ws.Range("N4").Value = [double]"3.004263504635674e-13"
$ws.Range("O4").Value = [double]"3.004263504635674e-13"
$ws.Range("U4").Value = "[44.05309089292267, 57.489691788979414]"

# Row 5
$ws.Range("M5").Value = "[52.685881833550106, 77.46862367115364]"
$ws.Range("N5").Value = [double]"8.681944052568724e-14"
$ws.Range("O5").Value = [double]"8.681944052568724e-14"
$ws.Range("Q5").Value = "[0.3333421634439624, 0.735868549489501]"
$ws.Range("R5").Value = [double]"2.845477669044527e-06"
$ws.Range("S5").Value = [double]"2.845477669044527e-06"
$ws.Range("U5").Value = "[44.833558309356775, 58.16402825204851]"
$ws.Range("Y5").Value = 22.41639639639687
$ws.Range("Z5").Value = 24.04298298298349

# Row 6
$ws.Range("M6").Value = "[52.69097293680791, 72.89668724556543]"
$ws.Range("Q6").Value = "[-0.09434212172942402, 0.25786846606042424]"
$ws.Range("R6").Value = 0.3547179786689276
$ws.Range("S6").Value = 0.3547179786689276
$ws.Range("U6").Value = "[41.0944860234605, 53.99158968587117]"
$ws.Range("Y6").Value = 22.10400400400416
$ws.Range("Z6").Value = 23.39609609609627

# Row 7
$ws.Range("M7").Value = "[52.28029439847279, 73.75210500867708]"
$ws.Range("N7").Value = [double]"1.998401444325282e-15"
$ws.Range("O7").Value = [double]"1.998401444325282e-15"
$ws.Range("U7").Value = "[43.325626862420535, 56.15517424500855]"

# Row 8
$ws.Range("M8").Value = "[51.48333470232198, 75.67251002443835]"
$ws.Range("N8").Value = [double]"8.43769498715119e-14"
$ws.Range("O8").Value = [double]"8.43769498715119e-14"
$ws.Range("Q8").Value = "[0.34592111300788364, 0.748447499053424]"
$ws.Range("R8").Value = [double]"1.859477410048882e-06"
$ws.Range("S8").Value = [double]"1.859477410048882e-06"
$ws.Range("U8").Value = "[43.64428865750558, 56.41864445022443]"
$ws.Range("Y8").Value = 20.30430430430445
$ws.Range("Z8").Value = 21.78098098098114
